# Add a note on iOS video playback to the feature list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 6 and 7 in column H: a short title and the detail note.
$ws.Range("H6").Value = "Video on iOS"
$ws.Range("H7").Value = "Can take video, can't immediately load and play. Decode error."

# Leave the cursor where the author left it after typing the note.
$ws.Range("H8").Select()
